# The commit re-sorts the weekly price records (rows 2-28) by re-shuffling
# which original record occupies which row, while keeping the constant
# columns (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria,
# Variedad, Calidad, Unidad de comercializacion, Kg o Unidades,
# Clasificacion) untouched, since they are identical for every row anyway.
#
# Only columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg) actually
# move between rows. Row 20 keeps its own original data.
#
# Mapping: target row -> source row (which original row's D/J/K/L/M/O/P
# values should end up in the target row).
$rowMap = @{
    2  = 24
    3  = 14
    4  = 16
    5  = 13
    6  = 17
    7  = 12
    8  = 22
    9  = 15
    10 = 2
    11 = 6
    12 = 27
    13 = 3
    14 = 21
    15 = 5
    16 = 4
    17 = 18
    18 = 26
    19 = 11
    20 = 20
    21 = 19
    22 = 8
    23 = 25
    24 = 7
    25 = 28
    26 = 23
    27 = 10
    28 = 9
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colD = 4
$colJ = 10
$colK = 11
$colL = 12
$colM = 13
$colO = 15
$colP = 16

# First snapshot the current (pre-edit) values of the columns that move,
# for every data row, before any writes happen - otherwise later writes
# would clobber values that still need to be read for other rows.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $row = @{}
    $row["D"] = $ws.Cells.Item($r, $colD).Value2
    $row["J"] = $ws.Cells.Item($r, $colJ).Value2
    $row["K"] = $ws.Cells.Item($r, $colK).Value2
    $row["L"] = $ws.Cells.Item($r, $colL).Value2
    $row["M"] = $ws.Cells.Item($r, $colM).Value2
    $row["O"] = $ws.Cells.Item($r, $colO).Value2
    $row["P"] = $ws.Cells.Item($r, $colP).Value2
    $snapshot[$r] = $row
}

# Now write the values back according to the mapping.
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $src = $snapshot[$sourceRow]

    $ws.Cells.Item($targetRow, $colD).Value2 = $src["D"]
    $ws.Cells.Item($targetRow, $colJ).Value2 = $src["J"]
    $ws.Cells.Item($targetRow, $colK).Value2 = $src["K"]
    $ws.Cells.Item($targetRow, $colL).Value2 = $src["L"]
    $ws.Cells.Item($targetRow, $colM).Value2 = $src["M"]
    $ws.Cells.Item($targetRow, $colO).Value2 = $src["O"]
    $ws.Cells.Item($targetRow, $colP).Value2 = $src["P"]
}
